# Nuevo formato 15 jun 2021
$wb = $excel.ActiveWorkbook

# --- Hoja "2o Parcial" (sheet2) ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

$ws2.Range("E3").Value = 36
$ws2.Range("F3").Value = 5
$ws2.Range("G3").Value = 87.8
$ws2.Range("H3").Value = 12.2
$ws2.Range("I3").Value = 8.300000000000001
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

# --- Hoja "3er Parcial" (sheet3) ---
$ws3 = $wb.Worksheets.Item("3er Parcial")

$ws3.Range("E2").Value = 38
$ws3.Range("F2").Value = 3
$ws3.Range("G2").Value = 92.68000000000001
$ws3.Range("H2").Value = 7.32
$ws3.Range("I2").Value = 7.3

$ws3.Range("E3").Value = 38
$ws3.Range("F3").Value = 3
$ws3.Range("G3").Value = 92.68000000000001
$ws3.Range("H3").Value = 7.32
$ws3.Range("I3").Value = 7.5

$ws3.Range("E4").Value = 27
$ws3.Range("F4").Value = 9
$ws3.Range("G4").Value = 75
$ws3.Range("H4").Value = 25
$ws3.Range("I4").Value = 6.9
